# Increase length of lines in RTS 24 data (TL_ESS sheet), again.
# For rows 40-87: the "length multiplier" used to derive MW ratings (col G/J)
# from the base circuit count (col M) changes from 1.75 to 2.0.
# Column H/F are formulas that recompute automatically from G.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TL_ESS")
$ws.Activate()

for ($r = 40; $r -le 87; $r++) {
    $m = $ws.Cells.Item($r, 13).Value2            # column M (base value, unchanged)
    $newG = $m * 2

    $ws.Cells.Item($r, 7).Value2 = $newG           # column G: hard-coded value = M*2
    $ws.Cells.Item($r, 10).Formula = "=M$r*2"      # column J: formula M{r}*2 (was M{r}*1.75)
}

# Restore the sheet view scroll position / selection recorded in the workbook.
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("L41").Select()
